# Regenerate save_data to use K (strikeouts) instead of Strike# in column G.
# Column G (header "K" in row 1) is updated row-by-row (rows 2..68) with
# the recalculated strikeout values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @(1,0,1,3,1,1,1,1,1,1,2,2,1,0,0,1,3,0,1,1,3,3,0,0,2,2,1,2,2,0,1,0,1,2,2,2,2,2,2,2,2,0,2,1,3,1,0,1,0,1,2,0,0,1,1,2,2,1,2,0,1,1,0,2,0,1,1)

$startRow = 2
for ($i = 0; $i -lt $newK.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 7).Value = $newK[$i]
}
